$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 197 (existing rows 197:226 shift down to 198:227)
$ws.Rows("197:197").Insert()

$ws.Range("A197").Value = 5
$ws.Range("B197").Value = "Macroferia Regional de Talca"
$ws.Range("C197").Value = "Maule"
$ws.Range("D197").Value = 44474
$ws.Range("E197").Value = 7
$ws.Range("F197").Value = 100112043
$ws.Range("G197").Value = "Pepino ensalada"
$ws.Range("H197").Value = "Sin especificar"
$ws.Range("I197").Value = "Primera"
$ws.Range("J197").Value = 300
$ws.Range("K197").Value = 17000
$ws.Range("L197").Value = 17000
$ws.Range("M197").Value = 17000
$ws.Range("N197").Value = "$/caja 60 unidades"
$ws.Range("O197").Value = "Región de Arica y Parinacota"
$ws.Range("P197").Value = 283
$ws.Range("Q197").Value = 60
$ws.Range("R197").Value = "Hortaliza"
